# Generate Report for Handoff
# - a new localized source file (d103f8d4-7c64-4fb3-9634-c8591153b075.md) failed its handoff
#   transform, so it gets a new row ("Handoff transform failed") inserted right after the
#   existing "Ready for handoff" row on every sheet.
# - the original "Ready for handoff" source file was re-generated with a new id
#   (2b17df3e-... -> 8576915b-...) and a new handoff .xlf + timestamp.

$wb = $excel.ActiveWorkbook

$oldUuid = "2b17df3e-2ca3-47a7-b48e-61097824f6ef"
$newUuid = "8576915b-9a0c-476b-9a4a-2d050f3c7a72"
$failedUuid = "d103f8d4-7c64-4fb3-9634-c8591153b075"

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/515af1e13d10a3d8a58c7e07eeb710e0511dd86f"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# drop every hyperlink up-front -- Range.Hyperlinks.Delete() clears the whole
# sheet's collection, which is fine since we rebuild all of them below anyway
$ws.Range("A1").Hyperlinks.Delete()

# insert the new row right after row 2 (pushes the old row 3 down to row 4,
# carrying its cell formatting along for free)
$ws.Rows.Item(3).Insert()

$ws.Range("A2").Value = "$newUuid.md"
$ws.Range("A3").Value = "$failedUuid.md"
$ws.Range("B3").Value = "Handoff transform failed"
$ws.Range("C3").Value = "Handoff transform failed"

$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/e2e/$newUuid.md", "", "", "$newUuid.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/e2e/$failedUuid.md", "", "", "$failedUuid.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "$repoBase/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A1").Hyperlinks.Delete()

$ws.Rows.Item(3).Insert()
$ws.Range("C3").Clear()   # Insert() leaves a stray empty styled cell behind here

$ws.Range("A2").Value = "$newUuid.md"
$ws.Range("C2").Value = "$newUuid.1c0d7312e88ffde0aacbd5c4f094f3e3a2202407.zh-cn.xlf"
$ws.Range("D2").Value = "2016-01-28 04:22:12"

$ws.Range("A3").Value = "$failedUuid.md"
$ws.Range("B3").Value = "Handoff transform failed"
$ws.Range("D3").Value = "0001-01-01 00:00:00"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Ignored"

$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ec163a579c386c556c75c69ee006804903f6789/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"

$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/e2e/$newUuid.md", "", "", "$newUuid.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "$zhHandoffBase/$newUuid.1c0d7312e88ffde0aacbd5c4f094f3e3a2202407.zh-cn.xlf", "", "", "$newUuid.1c0d7312e88ffde0aacbd5c4f094f3e3a2202407.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/e2e/$failedUuid.md", "", "", "$failedUuid.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "$repoBase/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A1").Hyperlinks.Delete()

$ws.Rows.Item(3).Insert()
$ws.Range("C3").Clear()

$ws.Range("A2").Value = "$newUuid.md"
$ws.Range("C2").Value = "$newUuid.1c0d7312e88ffde0aacbd5c4f094f3e3a2202407.de-de.xlf"
$ws.Range("D2").Value = "2016-01-28 04:22:22"

$ws.Range("A3").Value = "$failedUuid.md"
$ws.Range("B3").Value = "Handoff transform failed"
$ws.Range("D3").Value = "0001-01-01 00:00:00"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Ignored"

$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/51730c52d700829631e055982f1ce37154a41327/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang"

$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase/e2e/$newUuid.md", "", "", "$newUuid.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "$deHandoffBase/$newUuid.1c0d7312e88ffde0aacbd5c4f094f3e3a2202407.de-de.xlf", "", "", "$newUuid.1c0d7312e88ffde0aacbd5c4f094f3e3a2202407.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase/e2e/$failedUuid.md", "", "", "$failedUuid.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "$repoBase/.localization-config", "", "", ".localization-config")

$wb.Worksheets.Item("Overview").Activate()
